# Rename the "Data Segment Type ID" column to "Data Segment Type" on the
# AppNexus sheet, and switch the two sample rows from raw numeric type IDs
# (100 / 200) to their human-readable labels ("B2B" / "In-Market").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AppNexus")

$ws.Range("I1").Value = "Data Segment Type"
$ws.Range("I3").Value = "B2B"
$ws.Range("I4").Value = "In-Market"

$ws.Range("I5").Select()
